$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Vmax")
$ws2 = $wb.Worksheets.Item("Km")

# Add the new MANOVA row (row 10) to the Vmax sheet
$ws1.Range("A10").Value = "MANOVA"
$ws1.Range("B10").Value = "***"
$ws1.Range("C10").Value = "***"
$ws1.Range("D10").Value = "*"
$ws1.Range("E10").Value = "o"
$ws1.Range("F10").Value = "***"
$ws1.Range("G10").Value = "**"
$ws1.Range("H10").Value = "*"

# Add the same MANOVA row (row 10) to the Km sheet
$ws2.Range("A10").Value = "MANOVA"
$ws2.Range("B10").Value = "***"
$ws2.Range("C10").Value = "***"
$ws2.Range("D10").Value = "*"
$ws2.Range("E10").Value = "o"
$ws2.Range("F10").Value = "***"
$ws2.Range("G10").Value = "**"
$ws2.Range("H10").Value = "*"

# Set the selection on Vmax (no longer the active tab) to row 9
[void]$ws1.Rows.Item(9).Select()

# Activate Km and select row 9 there, making it the active tab/sheet
$ws2.Activate()
[void]$ws2.Rows.Item(9).Select()
